# Weekly update: insert this week's "Lechuga" (Macroferia Regional de Talca)
# price rows at the top of the data block (row 651), pushing the existing
# history down by 4 rows (710 -> 714).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 fresh rows at 651-654; everything previously at 651.. shifts to 655..
$ws.Rows("651:654").Insert()

# New week's data: 2021-09-22 (serial 44461)
$ws.Range("A651").Value = 5
$ws.Range("B651").Value = "Macroferia Regional de Talca"
$ws.Range("C651").Value = "Maule"
$ws.Range("D651").Value = 44461
$ws.Range("E651").Value = 7
$ws.Range("F651").Value = 100112033
$ws.Range("G651").Value = "Lechuga"
$ws.Range("H651").Value = "Conconina(o)"
$ws.Range("I651").Value = "Segunda"
$ws.Range("J651").Value = 500
$ws.Range("K651").Value = 4500
$ws.Range("L651").Value = 4500
$ws.Range("M651").Value = 4500
$ws.Range("N651").Value = "`$/caja 12 unidades"
$ws.Range("O651").Value = "Región del Maule"
$ws.Range("P651").Value = 375
$ws.Range("Q651").Value = 12
$ws.Range("R651").Value = "Hortaliza"

$ws.Range("A652").Value = 5
$ws.Range("B652").Value = "Macroferia Regional de Talca"
$ws.Range("C652").Value = "Maule"
$ws.Range("D652").Value = 44461
$ws.Range("E652").Value = 7
$ws.Range("F652").Value = 100112033
$ws.Range("G652").Value = "Lechuga"
$ws.Range("H652").Value = "Escarola"
$ws.Range("I652").Value = "Primera"
$ws.Range("J652").Value = 500
$ws.Range("K652").Value = 9000
$ws.Range("L652").Value = 9000
$ws.Range("M652").Value = 9000
$ws.Range("N652").Value = "`$/caja 15 unidades"
$ws.Range("O652").Value = "Provincia del Elquí"
$ws.Range("P652").Value = 600
$ws.Range("Q652").Value = 15
$ws.Range("R652").Value = "Hortaliza"

$ws.Range("A653").Value = 5
$ws.Range("B653").Value = "Macroferia Regional de Talca"
$ws.Range("C653").Value = "Maule"
$ws.Range("D653").Value = 44461
$ws.Range("E653").Value = 7
$ws.Range("F653").Value = 100112033
$ws.Range("G653").Value = "Lechuga"
$ws.Range("H653").Value = "Española"
$ws.Range("I653").Value = "Primera"
$ws.Range("J653").Value = 500
$ws.Range("K653").Value = 4500
$ws.Range("L653").Value = 4500
$ws.Range("M653").Value = 4500
$ws.Range("N653").Value = "`$/caja 18 unidades"
$ws.Range("O653").Value = "Región del Maule"
$ws.Range("P653").Value = 250
$ws.Range("Q653").Value = 18
$ws.Range("R653").Value = "Hortaliza"

$ws.Range("A654").Value = 5
$ws.Range("B654").Value = "Macroferia Regional de Talca"
$ws.Range("C654").Value = "Maule"
$ws.Range("D654").Value = 44461
$ws.Range("E654").Value = 7
$ws.Range("F654").Value = 100112033
$ws.Range("G654").Value = "Lechuga"
$ws.Range("H654").Value = "Marina"
$ws.Range("I654").Value = "Primera"
$ws.Range("J654").Value = 500
$ws.Range("K654").Value = 4500
$ws.Range("L654").Value = 4500
$ws.Range("M654").Value = 4500
$ws.Range("N654").Value = "`$/caja 18 unidades"
$ws.Range("O654").Value = "Región del Maule"
$ws.Range("P654").Value = 250
$ws.Range("Q654").Value = 18
$ws.Range("R654").Value = "Hortaliza"
